$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 updates
$ws.Range("F2").Value = 29
$ws.Range("G2").Value = 1000
$ws.Range("H2").Value = 1.11
$ws.Range("I2").Value = 1.13
$ws.Range("J2").Value = 1.04
$ws.Range("K2").Value = 840
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 1.08
$ws.Range("S2").Value = 1.08
$ws.Range("T2").Value = 0
$ws.Range("U2").Value = 0
$ws.Range("V2").Value = 1.01
$ws.Range("W2").Value = 1.01
$ws.Range("X2").Value = 1000
$ws.Range("Y2").Value = 1000
$ws.Range("Z2").Value = 1000
$ws.Range("AA2").Value = 1000
$ws.Range("AB2").Value = 1000
$ws.Range("AC2").Value = 1000
$ws.Range("AD2").Value = 1000
$ws.Range("AE2").Value = 1000
$ws.Range("AF2").Value = 1000
$ws.Range("AG2").Value = 1000
$ws.Range("AH2").Value = 1.08
$ws.Range("AI2").Value = 1000
$ws.Range("AJ2").Value = 1000
$ws.Range("AK2").Value = 1000
$ws.Range("AL2").Value = 1000
$ws.Range("AM2").Value = 1000
$ws.Range("AN2").Value = 1000
$ws.Range("AO2").Value = 1000

# Row 3 updates
$ws.Range("F3").Value = 1.57
$ws.Range("G3").Value = 1.6
$ws.Range("H3").Value = 12.5
$ws.Range("I3").Value = 14
$ws.Range("J3").Value = 3.3
$ws.Range("K3").Value = 3.5
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 1.27
$ws.Range("N3").Value = 1.72
$ws.Range("O3").Value = 2.36
$ws.Range("P3").Value = 1.2
$ws.Range("Q3").Value = 5.7
$ws.Range("R3").Value = 1.05
$ws.Range("S3").Value = 17.5
$ws.Range("T3").Value = 5.3
$ws.Range("U3").Value = 1.18
$ws.Range("V3").Value = 1.08
$ws.Range("W3").Value = 2.66
$ws.Range("X3").Value = 4.8
$ws.Range("Y3").Value = 21
$ws.Range("Z3").Value = 200
$ws.Range("AA3").Value = 1000
$ws.Range("AB3").Value = 3.3
$ws.Range("AC3").Value = 14.5
$ws.Range("AD3").Value = 990
$ws.Range("AE3").Value = 1000
$ws.Range("AF3").Value = 6.6
$ws.Range("AG3").Value = 990
$ws.Range("AH3").Value = 990
$ws.Range("AI3").Value = 1000
$ws.Range("AJ3").Value = 20
$ws.Range("AK3").Value = 70
$ws.Range("AL3").Value = 460
$ws.Range("AM3").Value = 1000
$ws.Range("AN3").Value = 40
$ws.Range("AO3").Value = 1000
